$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 562 (pushes the existing row 562..588 data
# down to 563..589, same as Excel's normal "Insert Row" behaviour).
$ws.Rows.Item(562).Insert()

# Populate the newly inserted row with the weekly price-report record.
$ws.Range("A562").Value = 3
$ws.Range("B562").Value = "Femacal de La Calera"
$ws.Range("C562").Value = "Coquimbo"
$ws.Range("D562").Value = 45147
$ws.Range("E562").Value = 5
$ws.Range("F562").Value = 100112012
$ws.Range("G562").Value = "Espinaca"
$ws.Range("H562").Value = "Sin especificar"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 65
$ws.Range("K562").Value = 4500
$ws.Range("L562").Value = 4500
$ws.Range("M562").Value = 4500
$ws.Range("N562").Value = "$/docena de atados (3 kilos)"
$ws.Range("O562").Value = "Provincia de Quillota"
$ws.Range("P562").Value = 1500
$ws.Range("Q562").Value = 3
$ws.Range("R562").Value = "Hortaliza"
